$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 168.66667
$ws.Range("I55").Value = 242
$ws.Range("J55").Value = 132
$ws.Range("K55").Value = 242
$ws.Range("L55").Value = 132
$ws.Range("M55").Value = -28
$ws.Range("N55").Value = -560

$ws.Range("H62").Value = 2900.9443
$ws.Range("I62").Value = 2703.7144
$ws.Range("J62").Value = 3591.25
$ws.Range("K62").Value = 2703.7144
$ws.Range("L62").Value = 3591.25
$ws.Range("M62").Value = -2079.7144
$ws.Range("N62").Value = -4839.25

$ws.Range("H65").Value = 2900.9443
$ws.Range("I65").Value = 2703.7144
$ws.Range("J65").Value = 3591.25
$ws.Range("K65").Value = 13518.572
$ws.Range("L65").Value = 17956.25
$ws.Range("M65").Value = -10398.572
$ws.Range("N65").Value = -24196.25

$ws.Range("H69").Value = 4910
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 4910
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 14730
$ws.Range("N69").Value = -16478
$ws.Range("M69").ClearContents()

$ws.Range("H72").Value = 4910
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 4910
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 44190
$ws.Range("N72").Value = -52926
$ws.Range("M72").ClearContents()

$ws.Range("H100").Value = 2335
$ws.Range("I100").Value = 1797.2727
$ws.Range("J100").Value = 4306.6665
$ws.Range("K100").Value = 1797.2727
$ws.Range("L100").Value = 4306.6665
$ws.Range("M100").Value = -1256.2727
$ws.Range("N100").Value = -5388.6665

$ws.Range("H113").Value = 3458.111
$ws.Range("I113").Value = 3459.9524
$ws.Range("J113").Value = 3451.6667
$ws.Range("K113").Value = 3459.9524
$ws.Range("L113").Value = 3451.6667
$ws.Range("M113").Value = -205.9524000000001
$ws.Range("N113").Value = -9959.6667

$ws.Range("H137").Value = 3363.9592
$ws.Range("I137").Value = 3668.1714
$ws.Range("K137").Value = 11004.5142
$ws.Range("M137").Value = -8454.514200000001

$ws.Range("H138").Value = 5126.1895
$ws.Range("I138").Value = 1476.2593
$ws.Range("J138").Value = 8305.161
$ws.Range("K138").Value = 4428.7779
$ws.Range("L138").Value = 24915.483
$ws.Range("M138").Value = 711.2221
$ws.Range("N138").Value = -35195.483

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3627.46
$ws.Range("I32").Value = 2909.2966
$ws.Range("J32").Value = 10888.889
$ws.Range("K32").Value = 2909.2966
$ws.Range("L32").Value = 10888.889
$ws.Range("M32").Value = -2622.2966
$ws.Range("N32").Value = -11462.889

$ws.Range("H61").Value = 1343.4359
$ws.Range("I61").Value = 696.97144
$ws.Range("J61").Value = 7000
$ws.Range("K61").Value = 696.97144
$ws.Range("L61").Value = 7000
$ws.Range("M61").Value = -484.97144
$ws.Range("N61").Value = -7424

$ws.Range("H132").Value = 2627.4146
$ws.Range("I132").Value = 1919.6
$ws.Range("J132").Value = 4557.8184
$ws.Range("K132").Value = 5758.799999999999
$ws.Range("L132").Value = 13673.4552
$ws.Range("M132").Value = -3228.799999999999
$ws.Range("N132").Value = -18733.4552

$ws.Range("H136").Value = 1343.4359
$ws.Range("I136").Value = 696.97144
$ws.Range("J136").Value = 7000
$ws.Range("K136").Value = 2090.91432
$ws.Range("L136").Value = 21000
$ws.Range("M136").Value = 459.0856800000001
$ws.Range("N136").Value = -26100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2106.2
$ws.Range("I107").Value = 1699.4897
$ws.Range("J107").Value = 3351.75
$ws.Range("K107").Value = 1699.4897
$ws.Range("L107").Value = 3351.75
$ws.Range("M107").Value = 220.5102999999999
$ws.Range("N107").Value = -7191.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3199.2917
$ws.Range("I31").Value = 2005.7587
$ws.Range("J31").Value = 5021
$ws.Range("K31").Value = 2005.7587
$ws.Range("L31").Value = 5021
$ws.Range("M31").Value = -1710.7587
$ws.Range("N31").Value = -5611

$ws.Range("H34").Value = 3199.2917
$ws.Range("I34").Value = 2005.7587
$ws.Range("J34").Value = 5021
$ws.Range("K34").Value = 2005.7587
$ws.Range("L34").Value = 5021
$ws.Range("M34").Value = -1803.7587
$ws.Range("N34").Value = -5425

$ws.Range("H107").Value = 1555.7826
$ws.Range("I107").Value = 423.125
$ws.Range("J107").Value = 4144.7144
$ws.Range("K107").Value = 423.125
$ws.Range("L107").Value = 4144.7144
$ws.Range("M107").Value = 1496.875
$ws.Range("N107").Value = -7984.7144

$ws.Range("H134").Value = 1395.9623
$ws.Range("I134").Value = 823.1429000000001
$ws.Range("J134").Value = 3583.0908
$ws.Range("K134").Value = 2469.4287
$ws.Range("L134").Value = 10749.2724
$ws.Range("M134").Value = 65.57129999999961
$ws.Range("N134").Value = -15819.2724

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 18162.9
$ws.Range("I120").Value = 14343.333
$ws.Range("K120").Value = 43029.999
$ws.Range("M120").Value = -38191.999

$ws.Range("H131").Value = 3144.1667
$ws.Range("J131").Value = 2620
$ws.Range("L131").Value = 7860
$ws.Range("N131").Value = -17940

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1215
$ws.Range("I107").Value = 525.25
$ws.Range("J107").Value = 1674.8334
$ws.Range("K107").Value = 525.25
$ws.Range("L107").Value = 1674.8334
$ws.Range("M107").Value = 1394.75
$ws.Range("N107").Value = -5514.8334

$ws.Range("H132").Value = 2615.0356
$ws.Range("I132").Value = 2319.0952
$ws.Range("J132").Value = 3502.8572
$ws.Range("K132").Value = 6957.285600000001
$ws.Range("L132").Value = 10508.5716
$ws.Range("M132").Value = -4427.285600000001
$ws.Range("N132").Value = -15568.5716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1375.5
$ws.Range("I7").Value = 886.6818
$ws.Range("J7").Value = 2450.9
$ws.Range("K7").Value = 886.6818
$ws.Range("L7").Value = 2450.9
$ws.Range("M7").Value = -774.6818
$ws.Range("N7").Value = -2674.9

$ws.Range("H126").Value = 1375.5
$ws.Range("I126").Value = 886.6818
$ws.Range("J126").Value = 2450.9
$ws.Range("K126").Value = 2660.0454
$ws.Range("L126").Value = 7352.700000000001
$ws.Range("M126").Value = -190.0454
$ws.Range("N126").Value = -12292.7

$ws.Range("H132").Value = 1470.4875
$ws.Range("I132").Value = 926.36505
$ws.Range("J132").Value = 3486.9412
$ws.Range("K132").Value = 2779.09515
$ws.Range("L132").Value = 10460.8236
$ws.Range("M132").Value = -249.0951500000001
$ws.Range("N132").Value = -15520.8236

$ws.Range("H136").Value = 1154.1299
$ws.Range("I136").Value = 783.82855
$ws.Range("J136").Value = 4857.143
$ws.Range("K136").Value = 2351.48565
$ws.Range("L136").Value = 14571.429
$ws.Range("M136").Value = 198.5143500000004
$ws.Range("N136").Value = -19671.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 8924.548000000001
$ws.Range("I132").Value = 1934.5103
$ws.Range("J132").Value = 23195.875
$ws.Range("K132").Value = 5803.5309
$ws.Range("L132").Value = 69587.625
$ws.Range("M132").Value = -3273.5309
$ws.Range("N132").Value = -74647.625

$ws.Range("H136").Value = 928.60376
$ws.Range("I136").Value = 416.56412
$ws.Range("K136").Value = 1249.69236
$ws.Range("M136").Value = 1300.30764

$ws.Range("H138").Value = 29750
$ws.Range("J138").Value = 29750
$ws.Range("L138").Value = 29750
$ws.Range("N138").Value = -40030

Write-Output "edits applied"
